$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2: D2, E2
$ws.Range("D2").Value = '27.607.43'
$ws.Range("E2").Value = '  -1.05%  '

# Row 3: D3, E3
$ws.Range("D3").Value = '1.664.59'
$ws.Range("E3").Value = '  -3.78%  '

# Row 4: E4
$ws.Range("E4").Value = '  +0.26%  '

# Row 5: D5, E5
Set-TextValue 5 4 '214.77'
$ws.Range("E5").Value = '  -1.25%  '

# Row 6: E6
$ws.Range("E6").Value = '  -1.83%  '

# Row 7: E7
$ws.Range("E7").Value = '  +0.22%  '

# Row 8: D8, E8
Set-TextValue 8 4 '23.41'
$ws.Range("E8").Value = '  -1.71%  '

# Row 9: D9, E9
Set-TextValue 9 4 '0.260'
$ws.Range("E9").Value = '  -2.09%  '

# Row 10: E10
$ws.Range("E10").Value = '  -1.49%  '

# Row 11: E11
$ws.Range("E11").Value = '  -2.60%  '

# Row 12: D12, E12
$ws.Range("D12").Value = '1.899.43'
$ws.Range("E12").Value = '  -3.77%  '

# Row 13: D13, E13
$ws.Range("D13").Value = '1.663.30'
$ws.Range("E13").Value = '  -3.82%  '

# Row 14: D14, E14
Set-TextValue 14 4 '4.13'
$ws.Range("E14").Value = '  -2.97%  '

# Row 15: E15
$ws.Range("E15").Value = '  -3.54%  '

# Row 16: D16, E16
Set-TextValue 16 4 '66.06'
$ws.Range("E16").Value = '  -2.99%  '

# Row 17: D17, E17
Set-TextValue 17 4 '247.15'
$ws.Range("E17").Value = '  +2.11%  '

# Row 18: D18, E18
$ws.Range("D18").Value = '27.601.30'
$ws.Range("E18").Value = '  -1.14%  '

# Row 19: E19
$ws.Range("E19").Value = '  -2.30%  '

# Row 20: D20, E20
Set-TextValue 20 4 '7.47'
$ws.Range("E20").Value = '  -6.91%  '

# Row 21: E21
$ws.Range("E21").Value = '  +0.18%  '

# Row 22: E22
$ws.Range("E22").Value = '  -3.76%  '

# Row 23: E23
$ws.Range("E23").Value = '  -3.98%  '

# Row 24: E24
$ws.Range("E24").Value = '  -4.74%  '

# Row 25: D25, E25
Set-TextValue 25 4 '146.33'
$ws.Range("E25").Value = '  -1.63%  '

# Row 26: E26
$ws.Range("E26").Value = '  -5.24%  '

# Row 27: D27, E27
Set-TextValue 27 4 '16.20'
$ws.Range("E27").Value = '  -2.73%  '

# Row 28: B28, C28, D28, E28
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 28 4 '0.112'
$ws.Range("E28").Value = '  -2.03%  '

# Row 29: B29, C29, D29, E29
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 29 4 '1.00'
$ws.Range("E29").Value = '  +0.23%  '

# Row 30: E30
$ws.Range("E30").Value = '  +4.34%  '

# Row 31: E31
$ws.Range("E31").Value = '  -1.27%  '

# Row 32: E32
$ws.Range("E32").Value = '  -3.41%  '

# Row 33: D33, E33
$ws.Range("D33").Value = '1.453.96'
$ws.Range("E33").Value = '  -2.24%  '

# Row 34: E34
$ws.Range("E34").Value = '  -5.21%  '

# Row 35: D35, E35
Set-TextValue 35 4 '1.57'
$ws.Range("E35").Value = '  -7.01%  '

# Row 36: D36, E36
Set-TextValue 36 4 '0.933'
$ws.Range("E36").Value = '  -3.58%  '

# Row 37: E37
$ws.Range("E37").Value = '  -0.98%  '

# Row 38: D38, E38
Set-TextValue 38 4 '0.580'
$ws.Range("E38").Value = '  -5.54%  '

# Row 39: D39, E39
Set-TextValue 39 4 '0.0169'
$ws.Range("E39").Value = '  -2.79%  '

# Row 40: E40
$ws.Range("E40").Value = '  -3.10%  '

# Row 41: D41, E41
Set-TextValue 41 4 '69.27'
$ws.Range("E41").Value = '  -3.96%  '

# Row 42: D42, E42
Set-TextValue 42 4 '0.999'
$ws.Range("E42").Value = '  +0.23%  '

# Row 43: D43, E43
Set-TextValue 43 4 '5.41'
$ws.Range("E43").Value = '  -7.84%  '

# Row 44: D44, E44
Set-TextValue 44 4 '0.794'
$ws.Range("E44").Value = '  +0.31%  '

# Row 45: D45, E45
$ws.Range("D45").Value = '1.807.47'
$ws.Range("E45").Value = '  -3.70%  '

# Row 46: E46
$ws.Range("E46").Value = '  -3.53%  '

# Row 47: E47
$ws.Range("E47").Value = '  +0.29%  '

# Row 48: D48, E48
Set-TextValue 48 4 '89.09'
$ws.Range("E48").Value = '  -2.99%  '

# Row 49: E49
$ws.Range("E49").Value = '  -3.31%  '

# Row 50: D50, E50
Set-TextValue 50 4 '0.102'
$ws.Range("E50").Value = '  -4.24%  '

# Row 51: D51, E51
Set-TextValue 51 4 '7.84'
$ws.Range("E51").Value = '  -5.86%  '
